$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price column (D) retains its original text formatting,
# since many of the new price strings (e.g. "1.000", "0.5140") look like
# numbers and would otherwise be auto-converted, dropping the formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.245.64'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.846.28'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '242.92'
$ws.Range("D6").Value = '0.6640'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '44.94'
$ws.Range("D9").Value = '0.07447'
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").Value = '0.2959'
$ws.Range("E10").Value = '  -0.16%  '
$ws.Range("D11").Value = '23.36'
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("D12").Value = '0.07777'
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").Value = '1.845.14'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = '5.026'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '0.6733'
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").Value = '83.52'
$ws.Range("E16").Value = '  -3.61%  '
$ws.Range("D17").Value = '6.194'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("D18").Value = '0.000008646'
$ws.Range("E18").Value = '  +4.73%  '
$ws.Range("D19").Value = '29.249.34'
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("D20").Value = '2.103.40'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '227.31'
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").Value = '7.192'
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("D25").Value = '1.000'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").Value = '158.91'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").Value = '0.1413'
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("D28").Value = '8.644'
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").Value = '4.141'
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("D32").Value = '4.061'
$ws.Range("D33").Value = '1.192'
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D34").Value = '0.05330'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").Value = '1.876'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").Value = '0.7481'
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("D37").Value = '1.159'
$ws.Range("E37").Value = '  +1.79%  '
$ws.Range("D38").Value = '2.655'
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("D39").Value = '1.321.52'
$ws.Range("E39").Value = '  -1.08%  '
$ws.Range("D40").Value = '0.01805'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '2.755'
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("D42").Value = '6.404'
$ws.Range("E42").Value = '  +6.58%  '
$ws.Range("D43").Value = '0.9031'
$ws.Range("E43").Value = '  -1.81%  '
$ws.Range("D44").Value = '0.9997'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '103.59'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").Value = '1.999.33'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = '65.56'
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("D48").Value = '0.00000000124'
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("B49").Value = 'XinFinNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D49").Value = '0.07755'
$ws.Range("E49").Value = '  -4.28%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.5140'
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").Value = '1.758'
$ws.Range("E51").Value = '  -0.62%  '
